# Add a new training block row "full_task_wo_driving_training" to the
# blocks table on Sheet1. This is inserted as the new row 6 (pushing the
# existing "full_task_training" row, and everything below it, down by one),
# and the word/nonword trial counts on the (now shifted) full_task_training
# row are reduced from 20 to 10 to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 6 ("full_task_training"),
# shifting rows 6:9 down to 7:10.
$ws.Rows(6).Insert()

# Populate the new row 6 with the full_task_wo_driving_training block.
# Same Georgia font / sound / instruction image / timing values as the
# (now-row-7) full_task_training block, but with fewer trials and its own
# task_name.
$ws.Range("A6").Value = "Georgia"
$ws.Range("B6").Value = "full_task_wo_driving_training"
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = "yes"
$ws.Range("H6").Value = "./instructions_pilot/driving_lexical_training.png"
$ws.Range("I6").Value = 1.77
$ws.Range("J6").Value = 2.07
$ws.Range("K6").Value = "no"
$ws.Range("L6").Value = $true

# The old row 6 (full_task_training) is now row 7; drop its trial counts
# from 20 to 10.
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 10

# Match the saved cursor position left behind by the edit.
[void]$ws.Range("E8").Select()
